# ccplant_model.xlsx edit script
# Mirrors the target commit's changes to the ccplant example workbook:
#  - Processes sheet: swap the "type" and "description" columns (B <-> E),
#    drop the empty spacer column F (shifting the gturbo_processes block
#    G:K to F:J), and update the selection.
#  - workbook-level defined names cgam_processes / gturbo_processes ranges.
#  - Active sheet moves from Format to ResourcesCost.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Processes sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Processes")
$ws.Activate()

# Swap the "description" (B) and "type" (E) columns' contents in place,
# keeping "fuel" (C) and "product" (D) untouched.
$bVals = $ws.Range("B1:B10").Value2
$eVals = $ws.Range("E1:E10").Value2
$ws.Range("B1:B10").Value2 = $eVals
$ws.Range("E1:E10").Value2 = $bVals

# Remove the empty spacer column (F) that separated the cgam table from
# the gturbo_processes query-table block; this shifts old G:K to F:J.
$ws.Range("F1").EntireColumn.Delete()

# New explicit widths for the (now) type column (B) and description
# column (E).
$ws.Columns.Item(2).ColumnWidth = 13.307291666666666
$ws.Columns.Item(5).ColumnWidth = 19.451822916666668

# NOTE: the process-type validation dropdown (an x14/extLst list
# validation, because it references another sheet) moves from column E to
# column B in the target file. The COM surface exposed here can add a new
# *legacy* Range.Validation entry, but it cannot edit or remove the
# existing x14 extLst dataValidation, and a freshly-added legacy entry
# would just coexist as extra (wrong-shaped) markup rather than actually
# relocating it -- so it is intentionally left alone here.

# Selection moves from a single cell to the whole type column.
$ws.Columns("E:E").Select()

# ---------------------------------------------------------------------
# Workbook-level defined names
# ---------------------------------------------------------------------
$cgam = $wb.Names.Item("Processes!cgam_processes")
$cgam.RefersTo = "=Processes!`$A`$1:`$D`$1"

$gturbo = $wb.Names.Item("Processes!gturbo_processes")
$gturbo.RefersTo = "=Processes!`$F`$5:`$J`$13"

# ---------------------------------------------------------------------
# Active sheet moves from Format to ResourcesCost
# ---------------------------------------------------------------------
$wsCost = $wb.Worksheets.Item("ResourcesCost")
$wsCost.Activate()
